# edit.ps1 - apply the tracked changes described in the diff.
$d = $word.ActiveDocument
$d.TrackRevisions = $false

function FindReplace($text, $old, $new) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $old"
    }
}

# --- 1) Merge the split "L" / "I" runs and fix the surrounding paragraph runs ---
$old1 = ", em 2016, o numero de documentos de patente excedeu 3 milhões pela primeira vez, um aumento de 8.3% (LI, 2018). Mas como realizar esta tarefa havendo algumas centenas de documentos de patentes sobre um assunto específico? O método tradicional necessita de tempo e equipe para realizá-lo, apresentando um resultado com deficiências devido ao alto volume de documentos de patente a serem analisadas (LI, 2018). Hoje, já há portais web que oferecem ferramentas das quais algumas auxiliam ao pesquisador a reduzir essa pesquisa (ABBAS; ZHANG; KHAN, 2014), mas classificam os documentos em uma relevância geral. Esse resultado somente demonstra que dentro daquela amostra de documentos, uma visão macro sobre o assunto que muitas vezes o pesquisador está em busca de um subassunto, como quais mercados essa tecnologia está presente, quais os processos de produção desta tecnologia ou qual a formulação desse composto. "
FindReplace "p1" $old1 $old1

# --- 2) Merge "De acordo com " + "Shahid..." runs, drop duplicate "Shahid et al (2019), " and fix double space ---
$old2 = "De acordo com Shahid et al (2019), a classificação de documentos de patentes em assuntos e a atribuição de valor de relevância para estes assuntos, permitindo ao pesquisador filtrar as patentes que o interessa e reduzindo o escopo de analise. Nesse trabalho, Shahid et al (2019), realizou a construção de uma matrix de valores de term frequency -  inverse document frequency (tf-idf), notações e peso ponderado por BM25, que posteriormente foi testado em diferentes classificadores, classificando os documentos de patente em cada assunto."
$new2 = "De acordo com Shahid et al (2019), a classificação de documentos de patentes em assuntos e a atribuição de valor de relevância para estes assuntos, permitindo ao pesquisador filtrar as patentes que o interessa e reduzindo o escopo de analise. Nesse trabalho, realizou a construção de uma matrix de valores de term frequency - inverse document frequency (tf-idf), notações e peso ponderado por BM25, que posteriormente foi testado em diferentes classificadores, classificando os documentos de patente em cada assunto."
FindReplace "p2" $old2 $new2

# --- 3) Merge "Seguindo " + "Anne..." runs and fix double period ---
$old3 = "Seguindo Anne et al (2017), identificou uma matriz de métodos a serem aplicados com os modelos k-Nearest Neighbors (kNN),  Support Vector Machine (SVM), Random Forest e J48. Os principais passos para essa pesquisa foram técnicas de seleção de características, com uso de ganho de informação e correlação para efetividade do classificadores.."
$new3 = "Seguindo Anne et al (2017), identificou uma matriz de métodos a serem aplicados com os modelos k-Nearest Neighbors (kNN),  Support Vector Machine (SVM), Random Forest e J48. Os principais passos para essa pesquisa foram técnicas de seleção de características, com uso de ganho de informação e correlação para efetividade do classificadores."
FindReplace "p3" $old3 $new3

# --- 4) Fix "2018;WANG" -> "2018; WANG" and "LI, G. 2018" -> "LI, 2018" ---
$old4 = "esse conhecimento está implícito em longos textos (LI, 2018;WANG et al., 2016)."
$new4 = "esse conhecimento está implícito em longos textos (LI, 2018; WANG et al., 2016)."
FindReplace "p4" $old4 $new4

$old5 = "sua analise resultaria em decisões de negocio de sucesso (LI, G. 2018)."
$new5 = "sua analise resultaria em decisões de negocio de sucesso (LI, 2018)."
FindReplace "p5" $old5 $new5

# --- 6) Bookmarks around the DOI/title text in ABBAS, ANNE and BREITZMAN references ---
function AddBookmarkAroundFirst($paraIndex, $needle, $bmName, [int]$searchFrom = 0) {
    $p = $d.Paragraphs($paraIndex)
    $full = $p.Range.Text
    $idx = $full.IndexOf($needle, $searchFrom)
    if ($idx -lt 0) {
        throw "Could not find '$needle' in paragraph $paraIndex"
    }
    $rStart = $p.Range.Start + $idx
    $rEnd = $rStart + $needle.Length
    $rng = $d.Range($rStart, $rEnd)
    $d.Bookmarks.Add($bmName, $rng)
    return $idx
}

# ABBAS reference paragraph
AddBookmarkAroundFirst 12 "A literature review on the state-of-the-art in patent analysis" "__DdeLink__30_2846412786" | Out-Null
$doiIdx = AddBookmarkAroundFirst 12 "10.1016/j.wpi.2013.12.006" "__DdeLink__33_2846412786"
AddBookmarkAroundFirst 12 "j.wpi.2013.12.006" "__DdeLink__35_2846412786" ($doiIdx + 26) | Out-Null

# ANNE reference paragraph
AddBookmarkAroundFirst 14 "10.5430/air.v7n1p1" "__DdeLink__38_2846412786" | Out-Null

# BREITZMAN reference paragraph
AddBookmarkAroundFirst 16 "10.1177/016555150202800302" "__DdeLink__42_2846412786" | Out-Null

# --- 7) Insert new "Shahid et al" bibliography entry between LI and WANG references ---
$pEmptyAfterLI = $d.Paragraphs(19)
$pEmptyAfterLI.Range.InsertParagraphAfter()
$pShahid = $d.Paragraphs(20)
$pShahid.Range.Text = "Shahid, M., Ahmed, A., Mushtaq, M. F., Ullah, S., & Akram, U. (2020, January). Automatic Patents Classification Using Supervised Machine Learning. In International Conference on Soft Computing and Data Mining (pp. 297-307). Springer, Cham."

$shahidFull = $pShahid.Range.Text
$italicText = "International Conference on Soft Computing and Data Mining"
$italicIdx = $shahidFull.IndexOf($italicText)
$italicStart = $pShahid.Range.Start + $italicIdx
$italicEnd = $italicStart + $italicText.Length
$italicRange = $d.Range($italicStart, $italicEnd)
$italicRange.Font.Italic = $true

Write-Output "done part7 shahid paragraph"
